$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row (2018) above the current first data row (2019),
# shifting all existing data rows down by one. We shift by copying full
# rows upward-to-downward starting from the bottom row so we don't
# overwrite data before it is copied, and we copy whole rows (not using
# Rows.Insert) so cell styles are preserved exactly and no new/duplicate
# style records get created.
$ws.Range("A7:C7").Copy($ws.Range("A8"))
$ws.Range("A6:E6").Copy($ws.Range("A7"))
$ws.Range("A5:E5").Copy($ws.Range("A6"))
$ws.Range("A4:E4").Copy($ws.Range("A5"))
$ws.Range("A3:E3").Copy($ws.Range("A4"))
$ws.Range("A2:E2").Copy($ws.Range("A3"))
$excel.CutCopyMode = 0

# The "*  revised" comment was anchored on C7; after the shift the cell it
# annotates is now C8, so relocate the comment accordingly.
$commentText = $ws.Range("C7").Comment.Text()
$ws.Range("C7").Comment.Delete()
$ws.Range("C8").AddComment($commentText)

# Fill in the new 2018 row of data.
$ws.Range("A2").Value = 2018
$ws.Range("B2").Value = 1.1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 2.8
$ws.Range("E2").Value = 2.1

# Update the selected cell recorded in the sheet view.
$ws.Range("I15").Select()

# Update the footer's generated-on timestamp.
$ws.PageSetup.LeftFooter = "Source: Bureau of Labor Statistics"
$ws.PageSetup.RightFooter = "Generated on: November 5, 2024 (12:50:38 AM)"
